$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy number formats/styles from column R into the new column S for rows 3-14
$ws.Range("R3:R14").Copy()
$ws.Range("S3").PasteSpecial(-4122)

# Fill in the new 2023 data (column S) for rows 3-14
$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 89.1
$ws.Range("S5").Value = 89.1
$ws.Range("S6").Value = 1895
$ws.Range("S7").Value = 1759
$ws.Range("S8").Value = 683.8
$ws.Range("S9").Value = 40.700000000000003
$ws.Range("S10").Value = 14.7
$ws.Range("S11").Value = 48.4
$ws.Range("S12").Value = 0.2
$ws.Range("S13").Value = 40.5
$ws.Range("S14").Value = "_"

# Update the current selection to match the new target cell
$ws.Range("D20").Select()
